# Update quarterly monthly numbers (report period: DEC 2019 -> Mar 2020)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TABLE A")

# --- Header / report-date text ---
$ws.Range("B2").Value = "Date of Report  April 1, 2020"
$ws.Range("A3").Value = "Navajo and Hopi Indian Relocation Program Status for Mar, 2020"

# --- Section 1: Completed Applications ---
# Certified (1)
$ws.Range("B5").Value = 3724
$ws.Range("E5").Value = 3842
$ws.Range("G5").Value = 1

# Denied
$ws.Range("B6").Value = 3210
$ws.Range("E6").Value = 3342
$ws.Range("F6").Value = "-"
$ws.Range("G6").Value = 1

# --- Section 2: Eligibility Appeals ---
# Active Cases
$ws.Range("B10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "-"
$ws.Range("G10").Value = 1

# --- Section 3: Status of Certified Applicants ---
# Certified
$ws.Range("B13").Value = 3724
$ws.Range("E13").Value = 3842
$ws.Range("G13").Value = 1

# Relocated
$ws.Range("B14").Value = -3593
$ws.Range("E14").Value = -3708
$ws.Range("G14").Value = 1

# Pending Relocation
$ws.Range("B16").Value = 5
$ws.Range("E16").Value = 5

# Contracts (2)
$ws.Range("B17").Value = -1
$ws.Range("E17").Value = -1
$ws.Range("F17").Value = "-"
$ws.Range("G17").Value = 1

# Remaining
$ws.Range("B19").Value = 3
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 1

# --- Section 4: Post Move Follow-up ---
# Relocated
$ws.Range("B22").Value = 3593
$ws.Range("E22").Value = 3708
$ws.Range("G22").Value = 1

# Closed Relocated
$ws.Range("B23").Value = -3579
$ws.Range("E23").Value = -3694
$ws.Range("G23").Value = 5

# Open PM
$ws.Range("B24").Value = 14
$ws.Range("E24").Value = 14
$ws.Range("G24").Value = 4

# Open 2+ Yrs (4)
$ws.Range("B25").Value = 3
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = "-"
$ws.Range("G25").Value = 4

# Denied Closed
$ws.Range("B26").Value = -3209
$ws.Range("E26").Value = -3341

# --- Restore the view: scroll back to top and select A4 ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("A4").Select() | Out-Null
